$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colNums = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$values = @{
    2 = @(1.02, 1.021690843904739, 1.024335558322903, 1.031020845341375, 1.038574876731988, 1.023594999628091, 1.026880887232523, 1.027164360691543, 1.033830170034838, 1.041362520405702, 1.012902122864229)
    3 = @(1.02, 1.022915699364301, 1.025472823016218, 1.032148780016427, 1.039878878316175, 1.023504579208683, 1.027742038119587, 1.028107796482225, 1.034765724350206, 1.042475259291091, 1.013198654409765)
    4 = @(1.02, 1.02370846351063, 1.02620923292253, 1.03287913268215, 1.040723474215049, 1.023443335729026, 1.028299001338133, 1.028718212256756, 1.035371015438443, 1.043195552127368, 1.013390147115894)
    5 = @(1.02, 1.024041791797002, 1.026518946372307, 1.033186294828393, 1.041078740419354, 1.023416932628374, 1.028533087849027, 1.028974820251848, 1.035625463062533, 1.043498430967106, 1.013470559150397)
    6 = @(1.02, 1.024097762092453, 1.02657095609456, 1.033237875890923, 1.041138402848366, 1.023412460900047, 1.028572388471618, 1.029017905266287, 1.035668184991953, 1.043549289674495, 1.013484055325069)
    7 = @(1.02, 1.023712917259372, 1.026213370829142, 1.032883236517173, 1.040728220518781, 1.023442985510381, 1.028302129451087, 1.028721641106224, 1.035374415444325, 1.04319959894115, 1.013391221944717)
    8 = @(1.02, 1.022104749281811, 1.024719794660998, 1.031401932699479, 1.039015401568393, 1.023565007130014, 1.027171971539207, 1.027483210386509, 1.034146361289477, 1.041738519439616, 1.013002416517993)
    9 = @(1.02, 1.019272392115521, 1.022091871198492, 1.028795482072916, 1.036003390439823, 1.023759139049229, 1.025178452617401, 1.025300496739729, 1.031981742428135, 1.039165957157512, 1.012314349007886)
    10 = @(1.02, 1.017384984856246, 1.020342488935981, 1.027060314229686, 1.03399943006943, 1.023874603872166, 1.023847995584439, 1.023844978841961, 1.030538160295124, 1.037452191705204, 1.011853645171917)
    11 = @(1.02, 1.016567882665679, 1.019585575667658, 1.026309532453585, 1.033132626194481, 1.023921309042974, 1.023271533123659, 1.023214616269054, 1.029912935901777, 1.036710392232467, 1.011653679462174)
    12 = @(1.02, 1.016264395459756, 1.019304509815653, 1.026030740622868, 1.032810793105168, 1.023938164268129, 1.023057352782322, 1.022980452815151, 1.02968067648427, 1.036434893748598, 1.011579331165185)
    13 = @(1.02, 1.016329493592705, 1.019364795558033, 1.026090538737243, 1.032879821258317, 1.023934571058996, 1.023103297775152, 1.023030682544738, 1.02973049799565, 1.03649398738241, 1.011595282396544)
    14 = @(1.02, 1.016542795907517, 1.019562340944287, 1.026286485772475, 1.033106020600884, 1.023922712353276, 1.023253830071777, 1.023195260636441, 1.02989373774941, 1.036687618660718, 1.011647535282064)
    15 = @(1.02, 1.016674221194607, 1.019684066443393, 1.026407226013614, 1.033245407492866, 1.023915340507297, 1.023346570452244, 1.023296660099796, 1.029994312020902, 1.036806926371481, 1.011679720441768)
    16 = @(1.02, 1.017439216323486, 1.020392734858089, 1.027110152733824, 1.034056976255929, 1.023871434968932, 1.023886245659489, 1.023886811393775, 1.030579651209672, 1.037501428055798, 1.011866906140499)
    17 = @(1.02, 1.017919117811125, 1.020837418041759, 1.027551228051581, 1.034566297647292, 1.023843013862359, 1.024224670625149, 1.02425696610334, 1.030946779363451, 1.037937142071138, 1.011984194602441)
    18 = @(1.02, 1.018199051623163, 1.021096850381806, 1.027808553720032, 1.034863465530954, 1.023826118601784, 1.024422032995243, 1.024472860222666, 1.031160905240654, 1.038191313094061, 1.012052560889342)
    19 = @(1.02, 1.018294504512702, 1.02118531968617, 1.027896304318693, 1.03496480732389, 1.023820303822167, 1.024489322541735, 1.024546472745361, 1.031233914345796, 1.038277983385228, 1.012075864218767)
    20 = @(1.02, 1.017867627352106, 1.020789701949646, 1.027503899283302, 1.034511643060462, 1.023846096030192, 1.024188364480832, 1.024217253142234, 1.030907391423385, 1.037890391364987, 1.011971615430321)
    21 = @(1.02, 1.016479983136562, 1.01950416640786, 1.026228782006593, 1.033039406784441, 1.02392621804881, 1.02320950363501, 1.023146797020599, 1.029845668365232, 1.036630598018782, 1.011632150103403)
    22 = @(1.02, 1.015607635870915, 1.018696392602322, 1.025427538414887, 1.032114541764025, 1.023973741382756, 1.022593727210534, 1.022473649605805, 1.02917798567144, 1.035838739815563, 1.011418297299844)
    23 = @(1.02, 1.016070073231402, 1.019124562618864, 1.025852248541195, 1.032604756279735, 1.023948818288725, 1.022920193533506, 1.022830508683172, 1.029531950268402, 1.036258498291528, 1.011531704399516)
    24 = @(1.02, 1.017890893626573, 1.020811262624287, 1.027525284949974, 1.03453633883968, 1.023844704313072, 1.024204769773426, 1.024235197753534, 1.030925189180637, 1.037911515915169, 1.011977299559188)
    25 = @(1.02, 1.020004467820774, 1.022770793549331, 1.029468871081095, 1.036781345557515, 1.02371141742819, 1.02569407411275, 1.025864840905052, 1.032541431143884, 1.039830793941406, 1.012492581078627)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $colNums.Length; $i++) {
        $ws.Cells.Item($row, $colNums[$i]).Value = $rowVals[$i]
    }
}
